# Add {herramientaSoftVersion} / {herramientaSoftDiscoVersion} placeholders
# right after the corresponding {herramientaSoft} / {herramientaSoftDisco}
# tokens in the five template paragraphs that describe:
#   - SIM card ("sim")
#   - the three loose-disk ("Un (01) {tipoDeDisco}, ...") blocks
#   - the storage-unit ("unidad de almacenamiento") block
#
# NOTE: there are look-alike paragraphs (the "disco suelto" block and the
# two extraction-summary sentences) that must stay untouched, so the
# replacements are scoped to specific paragraphs rather than done as a
# blanket document-wide Find/Replace.

$d = $word.ActiveDocument

$paraIndex = 0
foreach ($p in $d.Paragraphs) {
    $paraIndex = $paraIndex + 1
    $r = $p.Range

    if ($paraIndex -eq 5 -or $paraIndex -eq 16) {
        # "sim" block (5) and "unidad de almacenamiento" block (16):
        # {herramientaSoft}. -> {herramientaSoft} {herramientaSoftVersion}.
        $r.Find.Execute("herramientaSoft}.", $true, $false, $false, $false, $false, $true, 1, $false, "herramientaSoft} {herramientaSoftVersion}.", 2)
    }
    elseif ($paraIndex -eq 7 -or $paraIndex -eq 10 -or $paraIndex -eq 18) {
        # the three "Un (01) {tipoDeDisco}, ..." disk blocks:
        # {herramientaSoftDisco} -> {herramientaSoftDisco} {herramientaSoftDiscoVersion}
        $r.Find.Execute("herramientaSoftDisco}", $true, $false, $false, $false, $false, $true, 1, $false, "herramientaSoftDisco} {herramientaSoftDiscoVersion}", 2)
    }
}
